# "Origin story" slide (SlideID 262) - Content Placeholder 2:
# split "All from call with virgin media" into two runs reading
# "All came from " + "call with virgin media", and shrink the body
# text to 15pt (which PowerPoint turns into a normAutofit bodyPr).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Replace "All from " (chars 1-9) with "All came from " -> creates the
# first run; the untouched remainder "call with virgin media" becomes
# the second run automatically.
$tr.Characters(1, 9).Text = "All came from "

# Shrink both resulting runs to 15pt; PowerPoint auto-enables
# "Shrink text on overflow" (normAutofit) as a side effect.
$tr.Characters(1, 14).Font.Size = 15
$tr.Characters(15, 22).Font.Size = 15

Write-Output $tr.Text
